$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the previously-empty Train/Test accuracy cells for row 21
$ws.Range("D21").Value = 69.62
$ws.Range("E21").Value = 53.59
$ws.Range("F21").Value = 87.67
$ws.Range("G21").Value = 87.04

# Add a new row (22) documenting the modified scheduler.step() call
$ws.Range("A22").Value = 14
$ws.Range("B22").Value = "scheduler "
$ws.Range("B22").Font.ThemeColor = 1
$ws.Range("C22").Value = "scheduler.step(sum(losses)/len(losses))"
$ws.Range("C22").Font.ThemeColor = 1
$ws.Range("C22").HorizontalAlignment = -4131

$ws.Range("D22").Select()
